# Applies scheduled market-price refresh to the per-sheet Leve profit tables.
# For each touched row, currentAveragePrice* / LevePrice* / LeveProfit* (columns H:N)
# are refreshed from the latest market-board snapshot. Some rows gain or lose a
# LeveProfitHQ/NQ cell entirely when HQ or NQ pricing becomes (un)available.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 76
$ws.Range("H76").Value = 22387
$ws.Range("I76").Value = 14775
$ws.Range("J76").Value = 29999
$ws.Range("K76").Value = 14775
$ws.Range("L76").Value = 29999
$ws.Range("M76").Value = -14460
$ws.Range("N76").Value = -30629
# Row 79
$ws.Range("H79").Value = 22387
$ws.Range("I79").Value = 14775
$ws.Range("J79").Value = 29999
$ws.Range("K79").Value = 14775
$ws.Range("L79").Value = 29999
$ws.Range("M79").Value = -13683
$ws.Range("N79").Value = -32183
# Row 112
$ws.Range("H112").Value = 3853.611
$ws.Range("J112").Value = 3860.453
$ws.Range("L112").Value = 11581.359
$ws.Range("N112").Value = -13797.359
# Row 137
$ws.Range("H137").Value = 1950.2046
$ws.Range("I137").Value = 1688.3055
$ws.Range("J137").Value = 3128.75
$ws.Range("K137").Value = 5064.916499999999
$ws.Range("L137").Value = 9386.25
$ws.Range("M137").Value = -2514.916499999999
$ws.Range("N137").Value = -14486.25
# Row 138
$ws.Range("H138").Value = 3817.5615
$ws.Range("I138").Value = 1889.762
$ws.Range("J138").Value = 4942.1113
$ws.Range("K138").Value = 5669.286
$ws.Range("L138").Value = 14826.3339
$ws.Range("M138").Value = -529.2860000000001
$ws.Range("N138").Value = -25106.3339

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 774.3333
$ws.Range("I5").Value = 879.2
$ws.Range("J5").Value = 250
$ws.Range("K5").Value = 879.2
$ws.Range("L5").Value = 250
$ws.Range("M5").Value = -767.2
$ws.Range("N5").Value = -474
# Row 37
$ws.Range("H37").Value = 23949
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 23949
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 23949
$ws.Range("M37").ClearContents()
$ws.Range("N37").Value = -24495
# Row 74
$ws.Range("H74").Value = 1980.76
$ws.Range("I74").Value = 1502.25
$ws.Range("J74").Value = 3894.8
$ws.Range("K74").Value = 1502.25
$ws.Range("L74").Value = 3894.8
$ws.Range("M74").Value = -628.25
$ws.Range("N74").Value = -5642.8
# Row 77
$ws.Range("H77").Value = 1980.76
$ws.Range("I77").Value = 1502.25
$ws.Range("J77").Value = 3894.8
$ws.Range("K77").Value = 7511.25
$ws.Range("L77").Value = 19474
$ws.Range("M77").Value = -3143.25
$ws.Range("N77").Value = -28210
# Row 97
$ws.Range("H97").Value = 924.2143
$ws.Range("I97").Value = 929.8889
$ws.Range("J97").Value = 914
$ws.Range("K97").Value = 929.8889
$ws.Range("L97").Value = 914
$ws.Range("M97").Value = -433.8889
$ws.Range("N97").Value = -1906
# Row 98
$ws.Range("H98").Value = 12998
$ws.Range("J98").Value = 12998
$ws.Range("L98").Value = 12998
$ws.Range("N98").Value = -18988
# Row 102
$ws.Range("H102").Value = 14357.272
$ws.Range("I102").Value = 14357.272
$ws.Range("K102").Value = 14357.272
$ws.Range("M102").Value = -12735.272
# Row 132
$ws.Range("H132").Value = 4301.1
$ws.Range("I132").Value = 4001.2222
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 12003.6666
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -9473.6666
$ws.Range("N132").Value = -26060

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 774.3333
$ws.Range("I4").Value = 879.2
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 879.2
$ws.Range("L4").Value = 250
$ws.Range("M4").Value = -764.2
$ws.Range("N4").Value = -480
# Row 14
$ws.Range("H14").Value = 8080
$ws.Range("I14").Value = 8080
$ws.Range("K14").Value = 8080
$ws.Range("M14").Value = -7908
# Row 86
$ws.Range("H86").Value = 27779596
$ws.Range("I86").Value = 27779596
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 27779596
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -27778473
$ws.Range("N86").ClearContents()
# Row 89
$ws.Range("H89").Value = 27779596
$ws.Range("I89").Value = 27779596
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 138897980
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -138892364
$ws.Range("N89").ClearContents()
# Row 94
$ws.Range("H94").Value = 53572284
$ws.Range("I94").Value = 93750670
$ws.Range("J94").Value = 1103.1666
$ws.Range("K94").Value = 93750670
$ws.Range("L94").Value = 1103.1666
$ws.Range("M94").Value = -93750219
$ws.Range("N94").Value = -2005.1666
# Row 134
$ws.Range("H134").Value = 2768
$ws.Range("I134").Value = 2746.4443
$ws.Range("J134").Value = 2832.6667
$ws.Range("K134").Value = 8239.332900000001
$ws.Range("L134").Value = 8498.000100000001
$ws.Range("M134").Value = -5704.332900000001
$ws.Range("N134").Value = -13568.0001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 1039.6
$ws.Range("I7").Value = 1277
$ws.Range("K7").Value = 1277
$ws.Range("M7").Value = -1164
# Row 31
$ws.Range("H31").Value = 5304.8486
$ws.Range("I31").Value = 4071.1035
$ws.Range("J31").Value = 14249.5
$ws.Range("K31").Value = 4071.1035
$ws.Range("L31").Value = 14249.5
$ws.Range("M31").Value = -3776.1035
$ws.Range("N31").Value = -14839.5
# Row 34
$ws.Range("H34").Value = 5304.8486
$ws.Range("I34").Value = 4071.1035
$ws.Range("J34").Value = 14249.5
$ws.Range("K34").Value = 4071.1035
$ws.Range("L34").Value = 14249.5
$ws.Range("M34").Value = -3869.1035
$ws.Range("N34").Value = -14653.5
# Row 132
$ws.Range("H132").Value = 4498.5
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 4997
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 14991
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -20051
# Row 134
$ws.Range("H134").Value = 1425.1765
$ws.Range("I134").Value = 1493.8182
$ws.Range("J134").Value = 1299.3334
$ws.Range("K134").Value = 4481.4546
$ws.Range("L134").Value = 3898.0002
$ws.Range("M134").Value = -1946.4546
$ws.Range("N134").Value = -8968.0002

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 46
$ws.Range("H46").Value = 30413.285
$ws.Range("J46").Value = 61099.06
$ws.Range("L46").Value = 183297.18
$ws.Range("N46").Value = -183479.18
# Row 75
$ws.Range("H75").Value = 940
$ws.Range("J75").Value = 1266.6666
$ws.Range("L75").Value = 3799.9998
$ws.Range("N75").Value = -5795.9998
# Row 78
$ws.Range("H78").Value = 940
$ws.Range("J78").Value = 1266.6666
$ws.Range("L78").Value = 11399.9994
$ws.Range("N78").Value = -21383.9994
# Row 131
$ws.Range("H131").Value = 2506.15
$ws.Range("J131").Value = 2514.125
$ws.Range("L131").Value = 7542.375
$ws.Range("N131").Value = -17622.375

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1000
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 1000
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -1992
# Row 113
$ws.Range("H113").Value = 4072.7273
$ws.Range("I113").Value = 4492.5557
$ws.Range("J113").Value = 2183.5
$ws.Range("K113").Value = 4492.5557
$ws.Range("L113").Value = 2183.5
$ws.Range("M113").Value = -2322.5557
$ws.Range("N113").Value = -6523.5
# Row 126
$ws.Range("H126").Value = 1945.6666
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
# Row 132
$ws.Range("H132").Value = 6452.4
$ws.Range("I132").Value = 6452.4
$ws.Range("K132").Value = 19357.2
$ws.Range("M132").Value = -16827.2

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 12707.363
$ws.Range("I132").Value = 35760.332
$ws.Range("J132").Value = 4062.5
$ws.Range("K132").Value = 107280.996
$ws.Range("L132").Value = 12187.5
$ws.Range("M132").Value = -104750.996
$ws.Range("N132").Value = -17247.5
# Row 136
$ws.Range("H136").Value = 3632.2307
$ws.Range("I136").Value = 3753.8696
$ws.Range("K136").Value = 11261.6088
$ws.Range("M136").Value = -8711.6088

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 70
$ws.Range("H70").Value = 52991.5
$ws.Range("J70").Value = 52991.5
$ws.Range("L70").Value = 52991.5
$ws.Range("N70").Value = -53621.5
# Row 73
$ws.Range("H73").Value = 52991.5
$ws.Range("J73").Value = 52991.5
$ws.Range("L73").Value = 52991.5
$ws.Range("N73").Value = -55175.5
# Row 132
$ws.Range("H132").Value = 27999.6
$ws.Range("I132").Value = 27999.6
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 83998.79999999999
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -81468.79999999999
$ws.Range("N132").ClearContents()
# Row 133
$ws.Range("H133").Value = 95000
$ws.Range("J133").Value = 95000
$ws.Range("L133").Value = 95000
$ws.Range("N133").Value = -105120
